# Applies the 'Programa' (PT/EN) and 'Bibliografia' paragraph edits: each paragraph's
# single run of concatenated sentences is split into multiple <w:t> segments separated
# by manual line breaks (<w:br/>), matching the target diff exactly. The segment list
# for each paragraph is joined with "" to build the Find string (the original,
# unbroken text) and joined with Word's "^l" manual-line-break code to build the
# Replacement string, so Find/Replace performs the split in a single pass.

$d = $word.ActiveDocument

# --- "Programa" (Portuguese): 1./2./3./4./5. items, split by single manual line breaks ---
$programaPt = @(
    '1. Elementos de organizações de alto desempenho: aprendizagem organizacional, modelo da competência e capacitações dinâmicas.',
    '2. Introdução à Gestão de Pessoas',
    '3. Noções básicas de Marketing',
    '4. Desenvolvimento de atividade prática extensionista junto à micro e pequenos empreendedores da região (componente curricular: plano de marketing)',
    '5. Visita (viagem didática complementar) a uma empresa para conhecer e entender os diferentes processos organizacionais.'
)
$findPt = $programaPt -join ""
$replacePt = $programaPt -join "^l"
$rng = $d.Content
$okPt = $rng.Find.Execute($findPt, $true, $true, $false, $false, $false, $true, 1, $false, $replacePt, 2)
Write-Host "Programa (PT) split:" $okPt

# --- "Programa" (English, italic): 1./2./3./4./5. items, split by single manual line breaks ---
$programaEn = @(
    '1. Elements of high performance organizations: organizational learning, competency model and dynamic capabilities.',
    '2. Introduction to People Management',
    '3. Marketing basics',
    '4. Development of practical extension activities with micro and small entrepreneurs in the region (curricular component: marketing plan)',
    '5. Visit (complementary didactic trip) to a company to know and understand the different organizational processes.'
)
$findEn = $programaEn -join ""
$replaceEn = $programaEn -join "^l"
$rng = $d.Content
$okEn = $rng.Find.Execute($findEn, $true, $true, $false, $false, $false, $true, 1, $false, $replaceEn, 2)
Write-Host "Programa (EN) split:" $okEn

# --- "Bibliografia": each reference, split by a double manual line break (blank line) ---
$biblio = @(
    'Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.',
    'Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.',
    'ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.',
    'KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.',
    'KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.',
    'CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. ',
    'MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.',
    'GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.',
    'CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.',
    'SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.',
    'BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013',
    'KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005',
    'MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.'
)
$findBiblio = $biblio -join ""
$replaceBiblio = $biblio -join "^l^l"
$rng = $d.Content
$okBiblio = $rng.Find.Execute($findBiblio, $true, $true, $false, $false, $false, $true, 1, $false, $replaceBiblio, 2)
Write-Host "Bibliografia split:" $okBiblio

